# Generate Report for Handoff
# Refresh the localization-status report: a new handoff package (new GUID-named
# source file + new content hash) has just been generated. Update the
# Overview/zh-cn/de-de sheets accordingly and clear out the stale
# "latest target/handback" info for the two locales since the new handoff
# has not been handed back yet.

$wb = $excel.ActiveWorkbook

$oldGuid = "732581ef-e0f7-4d4e-b66e-fde12013220f"
$newGuid = "f7fc6c61-6621-45c0-8270-5b9d563016e5"
$newHash = "f13553d051c2d7d529b334a363f151ad197c120f"

$newFileName        = "$newGuid.md"
$newPathAndName      = "e2e\$newGuid.md"
$newZhHandoffFile    = "$newGuid.$newHash.zh-cn.xlf"
$newDeHandoffFile    = "$newGuid.$newHash.de-de.xlf"

$newGenerateDate      = "2016-11-29 05:04:32"
$zhHandoffDatetime    = "2016-11-29 05:04:17"
$zeroDatetime         = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName

$overviewUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/229fb3228dfd60093574706465e4ee72b6bdb355/e2e/$oldGuid.md"
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewUrl, "", "", $newPathAndName)

$wsOverview.Range("G2").Value = $newGenerateDate

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/caabcbf1d5c7d7438d08a8ae2890be1aafc51f0a/e2e/$oldGuid.md"
$overviewUrlZh = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/229fb3228dfd60093574706465e4ee72b6bdb355/e2e/$oldGuid.md"

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $overviewUrlZh, "", "", $newFileName)

$wsZh.Range("G2").Value = $newZhHandoffFile
$wsZh.Range("H2").Value = $zhHandoffDatetime

# Latest Target File / Latest Handback File reset to empty - new handoff has not
# come back from translation yet.
$wsZh.Range("I2").Hyperlinks.Delete()
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""

$wsZh.Range("K2").Value = $zeroDatetime
$wsZh.Range("O2").Value = "True"

$wsZh.Columns.Item(9).ColumnWidth = 18.65
$wsZh.Columns.Item(10).ColumnWidth = 21.7

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f83cae8c2f627c153f71785ede002a827d40874a/e2e/$oldGuid.md"
$overviewUrlDe = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/229fb3228dfd60093574706465e4ee72b6bdb355/e2e/$oldGuid.md"

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $overviewUrlDe, "", "", $newFileName)

$wsDe.Range("G2").Value = $newDeHandoffFile
$wsDe.Range("H2").Value = $newGenerateDate

$wsDe.Range("I2").Hyperlinks.Delete()
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""

$wsDe.Range("K2").Value = $zeroDatetime
$wsDe.Range("O2").Value = "True"

$wsDe.Columns.Item(9).ColumnWidth = 18.65
$wsDe.Columns.Item(10).ColumnWidth = 21.7

Write-Output "Report regenerated for handoff."
